# Auto-generated edit script: update FFXIV Leve market-price snapshot values
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets (scheduled market-data refresh).
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 365.84616
$ws.Range("J17").Value = 369.6842
$ws.Range("L17").Value = 1109.0526
$ws.Range("N17").Value = -1445.0526
$ws.Range("I18").Value = 1224.5
$ws.Range("K18").Value = 1224.5
$ws.Range("M18").Value = -940.5
$ws.Range("H70").Value = 1833
$ws.Range("I70").Value = 1600
$ws.Range("J70").Value = 1949.5
$ws.Range("K70").Value = 4800
$ws.Range("L70").Value = 5848.5
$ws.Range("M70").Value = -4530
$ws.Range("N70").Value = -6388.5
$ws.Range("H73").Value = 1833
$ws.Range("I73").Value = 1600
$ws.Range("J73").Value = 1949.5
$ws.Range("K73").Value = 4800
$ws.Range("L73").Value = 5848.5
$ws.Range("M73").Value = -3864
$ws.Range("N73").Value = -7720.5
$ws.Range("H132").Value = 5417
$ws.Range("I132").Value = 5343.1
$ws.Range("K132").Value = 16029.3
$ws.Range("M132").Value = -13499.3
$ws.Range("H135").Value = 624.8461
$ws.Range("I135").Value = 624.8461
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 5623.6149
$ws.Range("L135").Value = 0
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -3088.6149
$ws.Range("H137").Value = 4772.2793
$ws.Range("J137").Value = 7487.9565
$ws.Range("L137").Value = 22463.8695
$ws.Range("N137").Value = -27563.8695

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1803241.5
$ws.Range("I32").Value = 843894.7
$ws.Range("K32").Value = 843894.7
$ws.Range("M32").Value = -843607.7
$ws.Range("H61").Value = 1911.25
$ws.Range("I61").Value = 815
$ws.Range("K61").Value = 815
$ws.Range("M61").Value = -603
$ws.Range("H74").Value = 66675108
$ws.Range("I74").Value = 1235.1666
$ws.Range("K74").Value = 1235.1666
$ws.Range("M74").Value = -361.1666
$ws.Range("H77").Value = 66675108
$ws.Range("I77").Value = 1235.1666
$ws.Range("K77").Value = 6175.833000000001
$ws.Range("M77").Value = -1807.833000000001
$ws.Range("H132").Value = 2961485.8
$ws.Range("I132").Value = 4527590
$ws.Range("J132").Value = 3288.7778
$ws.Range("K132").Value = 13582770
$ws.Range("L132").Value = 9866.3334
$ws.Range("M132").Value = -13580240
$ws.Range("N132").Value = -14926.3334
$ws.Range("H136").Value = 1911.25
$ws.Range("I136").Value = 815
$ws.Range("K136").Value = 2445
$ws.Range("M136").Value = 105

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 41671870
$ws.Range("I20").Value = 49025344
$ws.Range("J20").Value = 2196.3333
$ws.Range("K20").Value = 49025344
$ws.Range("L20").Value = 2196.3333
$ws.Range("M20").Value = -49025097
$ws.Range("N20").Value = -2690.3333
$ws.Range("H82").Value = 14514.417
$ws.Range("I82").Value = 2736.3333
$ws.Range("J82").Value = 49848.668
$ws.Range("K82").Value = 2736.3333
$ws.Range("L82").Value = 49848.668
$ws.Range("M82").Value = -2353.3333
$ws.Range("N82").Value = -50614.668
$ws.Range("H85").Value = 14514.417
$ws.Range("I85").Value = 2736.3333
$ws.Range("J85").Value = 49848.668
$ws.Range("K85").Value = 2736.3333
$ws.Range("L85").Value = 49848.668
$ws.Range("M85").Value = -1410.3333
$ws.Range("N85").Value = -52500.668
$ws.Range("H134").Value = 2537.25
$ws.Range("I134").Value = 2537.25
$ws.Range("K134").Value = 7611.75
$ws.Range("M134").Value = -5076.75

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1763354.5
$ws.Range("I31").Value = 2118.6667
$ws.Range("K31").Value = 2118.6667
$ws.Range("M31").Value = -1823.6667
$ws.Range("H34").Value = 1763354.5
$ws.Range("I34").Value = 2118.6667
$ws.Range("K34").Value = 2118.6667
$ws.Range("M34").Value = -1916.6667
$ws.Range("H132").Value = 2414.6
$ws.Range("I132").Value = 2017.5128
$ws.Range("K132").Value = 6052.538399999999
$ws.Range("M132").Value = -3522.538399999999
$ws.Range("H134").Value = 3710.0967
$ws.Range("I134").Value = 3763.652
$ws.Range("J134").Value = 3556.125
$ws.Range("K134").Value = 11290.956
$ws.Range("L134").Value = 10668.375
$ws.Range("M134").Value = -8755.956
$ws.Range("N134").Value = -15738.375

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 106606910
$ws.Range("J32").Value = 19075112
$ws.Range("L32").Value = 57225336
$ws.Range("N32").Value = -57225902
$ws.Range("H55").Value = 5504.857
$ws.Range("J55").Value = 7029.778
$ws.Range("L55").Value = 21089.334
$ws.Range("N55").Value = -21443.334
$ws.Range("H68").Value = 1821460.6
$ws.Range("J68").Value = 2225857
$ws.Range("L68").Value = 6677571
$ws.Range("N68").Value = -6679193
$ws.Range("H71").Value = 1821460.6
$ws.Range("J71").Value = 2225857
$ws.Range("L71").Value = 20032713
$ws.Range("N71").Value = -20040825
$ws.Range("H122").Value = 965.9375
$ws.Range("J122").Value = 1074
$ws.Range("L122").Value = 9666
$ws.Range("N122").Value = -14566
$ws.Range("H134").Value = 3925
$ws.Range("I134").Value = 2696.4285
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 8089.2855
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -3019.2855
$ws.Range("N134").Value = -25140
$ws.Range("H139").Value = 10095.777

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 969.55554
$ws.Range("I31").Value = 840.75
$ws.Range("K31").Value = 840.75
$ws.Range("M31").Value = -548.75
$ws.Range("H37").Value = 969.55554
$ws.Range("I37").Value = 840.75
$ws.Range("K37").Value = 840.75
$ws.Range("M37").Value = -563.75
$ws.Range("H70").Value = 50015616
$ws.Range("I70").Value = 166669710
$ws.Range("J70").Value = 20999.857
$ws.Range("K70").Value = 166669710
$ws.Range("L70").Value = 20999.857
$ws.Range("M70").Value = -166669440
$ws.Range("N70").Value = -21539.857
$ws.Range("H73").Value = 50015616
$ws.Range("I73").Value = 166669710
$ws.Range("J73").Value = 20999.857
$ws.Range("K73").Value = 166669710
$ws.Range("L73").Value = 20999.857
$ws.Range("M73").Value = -166668774
$ws.Range("N73").Value = -22871.857
$ws.Range("H107").Value = 15326.857
$ws.Range("J107").Value = 18519.2
$ws.Range("L107").Value = 18519.2
$ws.Range("N107").Value = -22359.2
$ws.Range("H132").Value = 2050.8572
$ws.Range("I132").Value = 1730.238
$ws.Range("J132").Value = 3012.7144
$ws.Range("K132").Value = 5190.714
$ws.Range("L132").Value = 9038.143199999999
$ws.Range("M132").Value = -2660.714
$ws.Range("N132").Value = -14098.1432

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 133931310
$ws.Range("I22").Value = 23811138
$ws.Range("J22").Value = 200003420
$ws.Range("K22").Value = 23811138
$ws.Range("L22").Value = 200003420
$ws.Range("M22").Value = -23810843
$ws.Range("N22").Value = -200004010
$ws.Range("H27").Value = 133931310
$ws.Range("I27").Value = 23811138
$ws.Range("J27").Value = 200003420
$ws.Range("K27").Value = 23811138
$ws.Range("L27").Value = 200003420
$ws.Range("M27").Value = -23811031
$ws.Range("N27").Value = -200003634
$ws.Range("H46").Value = 6585.7144
$ws.Range("I46").Value = 6666.6665
$ws.Range("J46").Value = 6525
$ws.Range("K46").Value = 6666.6665
$ws.Range("L46").Value = 6525
$ws.Range("M46").Value = -6478.6665
$ws.Range("N46").Value = -6901
$ws.Range("H55").Value = 2814.6667
$ws.Range("I55").Value = 6927.6665
$ws.Range("J55").Value = 758.1667
$ws.Range("K55").Value = 6927.6665
$ws.Range("L55").Value = 758.1667
$ws.Range("M55").Value = -6754.6665
$ws.Range("N55").Value = -1104.1667
$ws.Range("H122").Value = 7712
$ws.Range("I122").Value = 6339
$ws.Range("J122").Value = 10000.333
$ws.Range("K122").Value = 19017
$ws.Range("L122").Value = 30000.999
$ws.Range("M122").Value = -16567
$ws.Range("N122").Value = -34900.999
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()
$ws.Range("H132").Value = 11757.48
$ws.Range("I132").Value = 10775.389
$ws.Range("K132").Value = 32326.167
$ws.Range("M132").Value = -29796.167
$ws.Range("H136").Value = 8427
$ws.Range("I136").Value = 6196.8
$ws.Range("K136").Value = 18590.4
$ws.Range("M136").Value = -16040.4

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5280.625
$ws.Range("I81").Value = 5115.6665
$ws.Range("J81").Value = 5555.5557
$ws.Range("K81").Value = 10231.333
$ws.Range("L81").Value = 11111.1114
$ws.Range("M81").Value = -9170.333000000001
$ws.Range("N81").Value = -13233.1114
$ws.Range("H84").Value = 5280.625
$ws.Range("I84").Value = 5115.6665
$ws.Range("J84").Value = 5555.5557
$ws.Range("K84").Value = 51156.665
$ws.Range("L84").Value = 55555.557
$ws.Range("M84").Value = -45852.665
$ws.Range("N84").Value = -66163.557
$ws.Range("H132").Value = 4712.4224
$ws.Range("I132").Value = 4105
$ws.Range("J132").Value = 6589.909
$ws.Range("K132").Value = 12315
$ws.Range("L132").Value = 19769.727
$ws.Range("M132").Value = -9785
$ws.Range("N132").Value = -24829.727
